$d = $word.ActiveDocument

# The trailing empty paragraph (just before the sectPr) becomes the
# first of three new Heading1 paragraphs: "Lied:", "Eingangsgebet", "Psalm".
$p = $d.Paragraphs.Last
$p.Range.Text = "Lied:"
$p.Style = "Heading1"

$p.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Eingangsgebet"
$p2.Style = "Heading1"

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "Psalm"
$p3.Style = "Heading1"
